# fitur baru: statistik produksi
# Update header labels with units, and add a new "Konsumsi Beras (Ton)" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename existing headers to include their units.
$ws.Range("C1").Value = "Luas Panen Padi (Ha)"
$ws.Range("D1").Value = "Produksi Padi (Ton GKG)"
$ws.Range("E1").Value = "Produksi Beras (Ton)"

# Add the new column header, copying E1's formatting first so the new
# cell shares the same style as the rest of the header row.
$ws.Range("E1").Copy($ws.Range("F1"))
$ws.Range("F1").Value = "Konsumsi Beras (Ton)"
